$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 41, shifting existing data (rows 41-131) down to rows 43-133.
$ws.Rows("41:42").Insert()

# Populate the first new row (41) with a new weekly data record.
$ws.Cells.Item(41,1).Value = 8
$ws.Cells.Item(41,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(41,3).Value = "Coquimbo"
$ws.Cells.Item(41,4).Value = 44519
$ws.Cells.Item(41,5).Value = 4
$ws.Cells.Item(41,6).Value = 100112021
$ws.Cells.Item(41,7).Value = "Ají"
$ws.Cells.Item(41,8).Value = "Inferno"
$ws.Cells.Item(41,9).Value = "Primera"
$ws.Cells.Item(41,10).Value = 500
$ws.Cells.Item(41,11).Value = 19000
$ws.Cells.Item(41,12).Value = 20000
$ws.Cells.Item(41,13).Value = 19500
$ws.Cells.Item(41,14).Value = "`$/caja 12 kilos"
$ws.Cells.Item(41,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(41,16).Value = 1625
$ws.Cells.Item(41,17).Value = 12
$ws.Cells.Item(41,18).Value = "Hortaliza"

# Populate the second new row (42) with a new weekly data record.
$ws.Cells.Item(42,1).Value = 8
$ws.Cells.Item(42,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(42,3).Value = "Coquimbo"
$ws.Cells.Item(42,4).Value = 44519
$ws.Cells.Item(42,5).Value = 4
$ws.Cells.Item(42,6).Value = 100112021
$ws.Cells.Item(42,7).Value = "Ají"
$ws.Cells.Item(42,8).Value = "Inferno"
$ws.Cells.Item(42,9).Value = "Segunda"
$ws.Cells.Item(42,10).Value = 360
$ws.Cells.Item(42,11).Value = 14000
$ws.Cells.Item(42,12).Value = 15000
$ws.Cells.Item(42,13).Value = 14500
$ws.Cells.Item(42,14).Value = "`$/caja 12 kilos"
$ws.Cells.Item(42,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(42,16).Value = 1208
$ws.Cells.Item(42,17).Value = 12
$ws.Cells.Item(42,18).Value = "Hortaliza"
